$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect, make the edits, then restore protection
$ws.Unprotect()

# Update the confidentiality disclaimer date from 2021-07-13 to 2021-07-14
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-14 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values
$ws.Range("D2").Value = 0.2479625592901645
$ws.Range("E2").Value = -0.002060855861316435

$ws.Range("D3").Value = 0.251544532487962
$ws.Range("E3").Value = -0.00463594218707386

$ws.Range("D4").Value = 0.2558868790943876
$ws.Range("E4").Value = -0.00009743739647261052

$ws.Range("D5").Value = 0.2446060291274859
$ws.Range("E5").Value = -0.02978804659155998

$ws.Range("E6").Value = -0.008988429747305515

# Restore sheet protection (objects/scenarios protected; column/row formatting allowed)
$ws.Protect($null, $true, $true, $true, $false, $false, $true, $true, $false, $false, $false, $false, $false, $false, $false, $false)
